$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert two new rows before the current row 14 ("resetbtcnt"), shifting the
# existing rows 14-16 down to 16-18.  Row insertion in this engine carries
# the formatting/values of the rows below along correctly.
# ---------------------------------------------------------------------------
$ws.Rows("14:15").Insert()

# ---------------------------------------------------------------------------
# New row 14: delay_time_num / u32 / "感应延时" / UPSSA0 / 0x30 / 32 / "秒"
# ---------------------------------------------------------------------------
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "delay_time_num"
$ws.Range("C14").Value = "u32"
$ws.Range("D14").Value = "感应延时"
$ws.Range("E14").Value = "UPSSA0"
$ws.Range("F14").Value = "0x30"
$ws.Range("G14").Value = 32
$ws.Range("H14").Value = "秒"

# ---------------------------------------------------------------------------
# New row 15: upload_duty / u32 / "通讯周期" / UPSSA0 / 0x34 / 8000 / "毫秒"
# ---------------------------------------------------------------------------
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "upload_duty"
$ws.Range("C15").Value = "u32"
$ws.Range("D15").Value = "通讯周期"
$ws.Range("E15").Value = "UPSSA0"
$ws.Range("F15").Value = "0x34"
$ws.Range("G15").Value = 8000
$ws.Range("H15").Value = "毫秒"

# ---------------------------------------------------------------------------
# Match the look of the surrounding table for the two new rows: thin border
# all around + left/center alignment for A:F (same visual style as the rest
# of the table, reusing the existing "bordered" cell style), and a plain
# left/center aligned style (no border) for the new H column, matching the
# rest of the sheet's H cells.
# ---------------------------------------------------------------------------
$ws.Range("A14:F15").Borders.LineStyle = 1
$ws.Range("A14:F15").HorizontalAlignment = -4131

$ws.Range("G14:G15").Borders.LineStyle = 1
$ws.Range("G14:G15").HorizontalAlignment = -4131

$ws.Range("H14:H15").HorizontalAlignment = -4131

# ---------------------------------------------------------------------------
# New column-I notes under the table (rows 19-22): source file / function
# references for the two new flash fields.
# ---------------------------------------------------------------------------
$ws.Range("I19").Value = "flash_def.h"
$ws.Range("I20").Value = "main.c:set_var_from_flash() and save_upssa0()"
$ws.Range("I21").Value = "protocol.c:reset_default_parameter()"
$ws.Range("I22").Value = "protocol.h"

$ws.Range("I19:I22").HorizontalAlignment = -4131

# ---------------------------------------------------------------------------
# Column widths - widen/introduce the new H column, keep the rest close to
# their previous sizing.
# ---------------------------------------------------------------------------
$ws.Columns("H:H").ColumnWidth = 9.977120535714286

# ---------------------------------------------------------------------------
# Selection cursor ends up on B16 (where "resetbtcnt" now lives) like in the
# authored workbook.
# ---------------------------------------------------------------------------
$ws.Range("B16").Select()
